# Apply updated dSF (column F) values as per repull of data / mean calculation fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -9
    8  = -2
    11 = -2
    13 = 6
    15 = -8
    16 = 0
    19 = -3
    20 = 4
    22 = -13
    23 = -5
    26 = -7
    31 = 1
    32 = -2
    35 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
